$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Work on the "general" sheet: insert two new rows (prior-distribution
#    settings) right above "Number of exp. conditions ..." (old row 6).
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("general")

$ws.Rows("6:7").Insert()

$ws.Range("A6").Value = "Prior distribution for fluxes (uniform or normal)"
$ws.Range("B6").Value = "normal"
$ws.Range("A7").Value = "Prior distribution for thermodynamic quantities (uniform or normal)"
$ws.Range("B7").Value = "normal"

$ws.Rows("6:7").RowHeight = 13.8

# Formatting: bold labels in column A (left/top aligned, bordered - matches
# the look of the other section headers on this sheet).
$labels = $ws.Range("A6:A7")
$labels.Font.Bold = $true
$labels.Font.Name = "Calibri"
$labels.Font.Size = 11
$labels.HorizontalAlignment = -4131
$labels.VerticalAlignment = -4160
$labels.Borders.LineStyle = 1

# Formatting: values in column B, centered, bordered.
$values = $ws.Range("B6:B7")
$values.Font.Name = "Calibri"
$values.Font.Size = 11
$values.HorizontalAlignment = -4108
$values.VerticalAlignment = -4107
$values.Borders.LineStyle = 1

# ---------------------------------------------------------------------------
# 2. Make the "general" sheet active again and re-select the newly edited
#    range (mirrors the end-user state captured in the saved workbook).
# ---------------------------------------------------------------------------
$ws.Activate()
$ws.Range("A6:B7").Select()
